# Updated cryptos list on Tue Aug 15 18:32:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "29.326.72";  E = "  -0.07%  " },
    @{ Row = 3;  D = "1.839.99";   E = "  -0.22%  " },
    @{ Row = 4;  D = "1.000";      E = "  +0.06%  " },
    @{ Row = 5;  D = "239.33";     E = "  -0.39%  " },
    @{ Row = 6;  D = "0.6221";     E = "  -0.84%  " },
    @{ Row = 7;  D = $null;        E = "  +0.05%  " },
    @{ Row = 8;  D = "0.07346";    E = "  -0.72%  " },
    @{ Row = 9;  D = "0.2876";     E = "  -0.51%  " },
    @{ Row = 10; D = $null;        E = "  -0.54%  " },
    @{ Row = 11; D = "0.07717";    E = "  -0.02%  " },
    @{ Row = 12; D = "1.822.13";   E = "  -1.25%  " },
    @{ Row = 13; D = "4.938";      E = $null },
    @{ Row = 14; D = $null;        E = "  +3.95%  " },
    @{ Row = 15; D = "0.6588";     E = "  -2.78%  " },
    @{ Row = 16; D = "81.34";      E = "  -0.90%  " },
    @{ Row = 17; D = "6.225";      E = "  -0.62%  " },
    @{ Row = 18; D = "29.316.24";  E = "  -0.20%  " },
    @{ Row = 19; D = "236.28";     E = "  +3.13%  " },
    @{ Row = 20; D = "12.18";      E = "  -0.99%  " },
    @{ Row = 21; D = $null;        E = "  +0.00%  " },
    @{ Row = 22; D = "7.193";      E = "  -3.51%  " },
    @{ Row = 23; D = "0.9979";     E = "  -0.33%  " },
    @{ Row = 24; D = "157.36";     E = "  -0.84%  " },
    @{ Row = 25; D = "8.398";      E = "  -0.87%  " },
    @{ Row = 26; D = "0.1330";     E = "  -1.61%  " },
    @{ Row = 27; D = "17.21";      E = "  -1.36%  " },
    @{ Row = 28; D = "0.06854";    E = $null },
    @{ Row = 29; D = $null;        E = "  +1.17%  " },
    @{ Row = 30; D = $null;        E = "  -0.52%  " },
    @{ Row = 31; D = "4.012";      E = "  -1.38%  " },
    @{ Row = 32; D = "3.929";      E = "  -3.49%  " },
    @{ Row = 33; D = "1.151";      E = "  +1.04%  " },
    @{ Row = 34; D = "1.742";      E = "  -5.18%  " },
    @{ Row = 35; D = $null;        E = "  -2.08%  " },
    @{ Row = 36; D = $null;        E = "  +0.23%  " },
    @{ Row = 37; D = "0.01816";    E = "  -2.05%  " },
    @{ Row = 38; D = "2.781";      E = "  -1.36%  " },
    @{ Row = 39; D = "1.231.26";   E = "  -0.98%  " },
    @{ Row = 40; D = "6.642";      E = $null },
    @{ Row = 41; D = "0.9426";     E = "  +0.71%  " },
    @{ Row = 42; D = $null;        E = "  +0.19%  " },
    @{ Row = 43; D = "1.988.27";   E = "  -1.53%  " },
    @{ Row = 44; D = "101.25";     E = "  +0.71%  " },
    @{ Row = 45; D = $null;        E = "  -1.14%  " },
    @{ Row = 46; D = $null;        E = "  +4.68%  " },
    @{ Row = 47; D = "1.684";      E = "  -1.82%  " },
    @{ Row = 48; D = "6.862";      E = "  -2.49%  " },
    @{ Row = 49; D = "8.811";      E = "  -2.26%  " },
    @{ Row = 50; D = "0.1125";     E = "  -2.14%  " },
    @{ Row = 51; D = "0.3848";     E = "  -1.55%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
